# Update cryptos list figures (Price / Volume(1h)) for the Fri Mar  1 12:53:12 UTC 2024 run.
# Numeric-looking Price strings are prefixed with an apostrophe so Excel keeps
# them as text (matching the source data, which stores Price as plain text
# even for values like "407.64"), then the auto-applied Text number format is
# reset back to Normal so no extra cell styling is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.503.64"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "3.433.76"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'407.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").Value = "'134.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.46%  "
$ws.Range("D7").Value = "'0.592"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.689"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("E10").Value = "  -3.44%  "
$ws.Range("D11").Value = "'42.33"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.89%  "
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("D13").Value = "'8.48"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.42%  "
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").Value = "3.433.40"
$ws.Range("E15").Value = "  -2.19%  "
$ws.Range("D16").Value = "62.419.45"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").Value = "'11.38"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.89%  "
$ws.Range("E18").Value = "  -2.11%  "
$ws.Range("D19").Value = "'0.0000132"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.77%  "
$ws.Range("E20").Value = "  -4.60%  "
$ws.Range("D21").Value = "'84.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.72%  "
$ws.Range("D22").Value = "'314.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.69%  "
$ws.Range("D23").Value = "'12.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.21%  "
$ws.Range("D24").Value = "'3.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.48%  "
$ws.Range("E25").Value = "  +8.99%  "
$ws.Range("D26").Value = "'29.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.66%  "
$ws.Range("D27").Value = "'8.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.29%  "
$ws.Range("D28").Value = "'2.74"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.51%  "
$ws.Range("D29").Value = "'7.57"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.94%  "
$ws.Range("E30").Value = "  -2.94%  "
$ws.Range("E31").Value = "  -2.68%  "
$ws.Range("D32").Value = "'42.54"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.17%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").Value = "'11.38"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.39%  "
$ws.Range("E35").Value = "  -1.31%  "
$ws.Range("D36").Value = "'51.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.94%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").Value = "'3.41"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.38%  "
$ws.Range("D39").Value = "'2.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.42%  "
$ws.Range("D40").Value = "'0.313"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.18%  "
$ws.Range("D41").Value = "'138.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.59%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D44").Value = "'4.05"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.43%  "
$ws.Range("D45").Value = "'16.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.49%  "
$ws.Range("E46").Value = "  -0.99%  "
$ws.Range("D47").Value = "'21.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.60%  "
$ws.Range("D48").Value = "2.123.00"
$ws.Range("D49").Value = "'2.31"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.17%  "
$ws.Range("D50").Value = "'1.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.17%  "
$ws.Range("D51").Value = "'1.68"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +18.65%  "
